# This edit reorders the species-observation records (rows 2-41) on the
# "Artfynd" sheet: the same 40 records are redistributed across the rows,
# i.e. the whole sheet body is permuted row-wise (row 27 happens to stay
# where it is). The mapping below gives, for every row in the NEW layout,
# which row in the CURRENT (pre-edit) layout supplies its data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    2=19;  3=11;  4=15;  5=12;  6=8;   7=5;   8=41;  9=25;  10=3;  11=23;
    12=22; 13=17; 14=31; 15=33; 16=13; 17=2;  18=4;  19=26; 20=18; 21=6;
    22=30; 23=29; 24=35; 25=37; 26=32; 27=27; 28=9;  29=28; 30=39; 31=10;
    32=20; 33=7;  34=16; 35=40; 36=34; 37=24; 38=21; 39=36; 40=38; 41=14
}

$firstRow = 2
$lastRow = 41

# Snapshot all current data (A1:AY41) before writing anything, so the
# reshuffle can freely read from rows we are about to overwrite.
$fullRange = $ws.Range("A1:AY41")
$source = $fullRange.Value()

# The data is written back in three separate blocks that skip columns Y
# and AA (Startdatum / Slutdatum). Those two columns hold plain text dates
# ("2023-08-16") and, if they are included in a bulk Range.Value array
# write, Excel auto-converts the look-alike text into a real date value -
# changing the cell's stored type even though the text itself would be
# identical either way. Since every row already shares the exact same
# Startdatum/Slutdatum text, those two columns never actually need to be
# touched by the permutation, so simply leaving them out of the write
# keeps them pixel/byte perfect.

# --- Block 1: columns A..X (1-24) ---
$range1 = $ws.Range("A1:X41")
$block1 = $range1.Value()
for ($newRow = $firstRow; $newRow -le $lastRow; $newRow++) {
    $oldRow = $rowMap[$newRow]
    for ($c = 1; $c -le 24; $c++) {
        $block1[$newRow, $c] = $source[$oldRow, $c]
    }
}
$range1.Value = $block1

# --- Block 2: column Z (26) ---
$range2 = $ws.Range("Z1:Z41")
$block2 = $range2.Value()
for ($newRow = $firstRow; $newRow -le $lastRow; $newRow++) {
    $oldRow = $rowMap[$newRow]
    $block2[$newRow, 1] = $source[$oldRow, 26]
}
$range2.Value = $block2

# --- Block 3: columns AB..AY (28-51) ---
$range3 = $ws.Range("AB1:AY41")
$block3 = $range3.Value()
for ($newRow = $firstRow; $newRow -le $lastRow; $newRow++) {
    $oldRow = $rowMap[$newRow]
    for ($c = 28; $c -le 51; $c++) {
        $block3[$newRow, $c - 27] = $source[$oldRow, $c]
    }
}
$range3.Value = $block3
